$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "M2"  = "[49.84255592331836, 50.08151485539584]"
    "U2"  = "[49.922024497018036, 50.09031183432045]"
    "M3"  = "[49.875664421903, 50.07190388909024]"
    "U3"  = "[49.93310901406607, 50.06735964803976]"
    "M4"  = "[49.84815503663116, 50.11730889437312]"
    "U4"  = "[49.90745821332454, 50.083747913721545]"
    "M5"  = "[49.863437809908945, 50.132380004074804]"
    "U5"  = "[49.905415160065935, 50.07539421158362]"
    "M6"  = "[49.86346195963935, 50.08887281641991]"
    "U6"  = "[49.91429305903233, 50.05479923007705]"
    "M7"  = "[49.813372116277925, 50.09485966392766]"
    "U7"  = "[49.89530848932748, 50.06477522912849]"
    "M8"  = "[49.91017641595767, 50.21489357268637]"
    "U8"  = "[49.906433922457744, 50.08497709381542]"
    "M9"  = "[49.88282377623477, 50.1681435606814]"
    "U9"  = "[49.908197333595965, 50.071933757304436]"
    "M10" = "[49.904133902890706, 50.210359373295546]"
    "U10" = "[49.8778570178721, 50.05043659315638]"
    "M11" = "[49.945908435898794, 50.216678310779635]"
    "U11" = "[49.85110373343791, 50.009292324662155]"
    "M12" = "[49.9212890774517, 50.23692030536118]"
    "U12" = "[49.934717231320136, 50.111344119953245]"
    "M13" = "[49.848070252155644, 50.14002283570259]"
    "U13" = "[49.953130445935464, 50.112571532080146]"
    "M14" = "[49.88296223923192, 50.18676315962532]"
    "U14" = "[49.919139607286134, 50.07927112304129]"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
